$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Release_EDI")

# Row 1
$ws.Cells.Item(1,1).Value = "projectDescriptionID"
$ws.Cells.Item(1,2).Value = "releaseID"
$ws.Cells.Item(1,3).Value = "commonName"
$ws.Cells.Item(1,4).Value = "markedRun"
$ws.Cells.Item(1,5).Value = "markedLifeStage"
$ws.Cells.Item(1,6).Value = "markedFishOrigin"
$ws.Cells.Item(1,7).Value = "sourceOfFishSite"
$ws.Cells.Item(1,8).Value = "releaseSite"
$ws.Cells.Item(1,9).Value = "releaseSubSite"
$ws.Cells.Item(1,10).Value = "nReleased"
$ws.Cells.Item(1,11).Value = "releaseTime"
$ws.Cells.Item(1,12).Value = "testDays"
$ws.Cells.Item(1,13).Value = "appliedMarkType"
$ws.Cells.Item(1,14).Value = "appliedMarkColor"
$ws.Cells.Item(1,15).Value = "appliedMarkPosition"
$ws.Cells.Item(1,16).Value = "appliedMarkCode"
$ws.Cells.Item(1,17).Value = "includeAnalysis"

# Row 2
$ws.Cells.Item(2,1).Value = 11
$ws.Cells.Item(2,2).Value = 256
$ws.Cells.Item(2,3).Value = "Chinook salmon"
$ws.Cells.Item(2,4).Value = "Spring"
$ws.Cells.Item(2,5).Value = "Juvenile"
$ws.Cells.Item(2,6).Value = "Natural"
$ws.Cells.Item(2,7).Value = "Parrott-Phelan canal trap box"
$ws.Cells.Item(2,8).Value = "Not applicable"
$ws.Cells.Item(2,9).Value = $null
$ws.Cells.Item(2,10).Value = 100
$ws.Cells.Item(2,11).Value = 44202.5750578704
$ws.Cells.Item(2,12).Value = 7
$ws.Cells.Item(2,13).Value = "Pigment / dye"
$ws.Cells.Item(2,14).Value = "Brown"
$ws.Cells.Item(2,15).Value = "Whole body"
$ws.Cells.Item(2,16).Value = $null
$ws.Cells.Item(2,17).Value = "Yes"

# Row 3
$ws.Cells.Item(3,1).Value = 11
$ws.Cells.Item(3,2).Value = 257
$ws.Cells.Item(3,3).Value = "Chinook salmon"
$ws.Cells.Item(3,4).Value = "Spring"
$ws.Cells.Item(3,5).Value = $null
$ws.Cells.Item(3,6).Value = "Natural"
$ws.Cells.Item(3,7).Value = $null
$ws.Cells.Item(3,8).Value = "Not applicable"
$ws.Cells.Item(3,9).Value = $null
$ws.Cells.Item(3,10).Value = 100
$ws.Cells.Item(3,11).Value = 44202.5208680556
$ws.Cells.Item(3,12).Value = 7
$ws.Cells.Item(3,13).Value = "Pigment / dye"
$ws.Cells.Item(3,14).Value = "Brown"
$ws.Cells.Item(3,15).Value = "Whole body"
$ws.Cells.Item(3,16).Value = $null
$ws.Cells.Item(3,17).Value = "Yes"

# Row 4
$ws.Cells.Item(4,1).Value = 11
$ws.Cells.Item(4,2).Value = 258
$ws.Cells.Item(4,3).Value = "Chinook salmon"
$ws.Cells.Item(4,4).Value = "Spring"
$ws.Cells.Item(4,5).Value = $null
$ws.Cells.Item(4,6).Value = "Natural"
$ws.Cells.Item(4,7).Value = $null
$ws.Cells.Item(4,8).Value = "Not applicable"
$ws.Cells.Item(4,9).Value = $null
$ws.Cells.Item(4,10).Value = 100
$ws.Cells.Item(4,11).Value = 44481.5832060185
$ws.Cells.Item(4,12).Value = 7
$ws.Cells.Item(4,13).Value = "Pigment / dye"
$ws.Cells.Item(4,14).Value = "Brown"
$ws.Cells.Item(4,15).Value = "Whole body"
$ws.Cells.Item(4,16).Value = $null
$ws.Cells.Item(4,17).Value = "No"

# Row 5
$ws.Cells.Item(5,1).Value = 11
$ws.Cells.Item(5,2).Value = 259
$ws.Cells.Item(5,3).Value = "Not applicable (n/a)"
$ws.Cells.Item(5,4).Value = $null
$ws.Cells.Item(5,5).Value = $null
$ws.Cells.Item(5,6).Value = $null
$ws.Cells.Item(5,7).Value = $null
$ws.Cells.Item(5,8).Value = "Not applicable"
$ws.Cells.Item(5,9).Value = $null
$ws.Cells.Item(5,10).Value = $null
$ws.Cells.Item(5,11).Value = $null
$ws.Cells.Item(5,12).Value = $null
$ws.Cells.Item(5,13).Value = $null
$ws.Cells.Item(5,14).Value = $null
$ws.Cells.Item(5,15).Value = $null
$ws.Cells.Item(5,16).Value = $null
$ws.Cells.Item(5,17).Value = $null

# Row 6
$ws.Cells.Item(6,1).Value = 11
$ws.Cells.Item(6,2).Value = 260
$ws.Cells.Item(6,3).Value = "Not applicable (n/a)"
$ws.Cells.Item(6,4).Value = $null
$ws.Cells.Item(6,5).Value = $null
$ws.Cells.Item(6,6).Value = $null
$ws.Cells.Item(6,7).Value = $null
$ws.Cells.Item(6,8).Value = "Not applicable"
$ws.Cells.Item(6,9).Value = $null
$ws.Cells.Item(6,10).Value = $null
$ws.Cells.Item(6,11).Value = $null
$ws.Cells.Item(6,12).Value = $null
$ws.Cells.Item(6,13).Value = $null
$ws.Cells.Item(6,14).Value = $null
$ws.Cells.Item(6,15).Value = $null
$ws.Cells.Item(6,16).Value = $null
$ws.Cells.Item(6,17).Value = $null

# Row 7
$ws.Cells.Item(7,1).Value = 11
$ws.Cells.Item(7,2).Value = 261
$ws.Cells.Item(7,3).Value = "Chinook salmon"
$ws.Cells.Item(7,4).Value = "Spring"
$ws.Cells.Item(7,5).Value = "Juvenile"
$ws.Cells.Item(7,6).Value = "Natural"
$ws.Cells.Item(7,7).Value = "Parrot-Phelan RST"
$ws.Cells.Item(7,8).Value = "Parrott-Phelan e-test release site"
$ws.Cells.Item(7,9).Value = $null
$ws.Cells.Item(7,10).Value = 100
$ws.Cells.Item(7,11).Value = 44202.5208680556
$ws.Cells.Item(7,12).Value = 6
$ws.Cells.Item(7,13).Value = "Pigment / dye"
$ws.Cells.Item(7,14).Value = "Brown"
$ws.Cells.Item(7,15).Value = "Whole body"
$ws.Cells.Item(7,16).Value = $null
$ws.Cells.Item(7,17).Value = "Yes"

# Row 8
$ws.Cells.Item(8,1).Value = 11
$ws.Cells.Item(8,2).Value = 262
$ws.Cells.Item(8,3).Value = "Not applicable (n/a)"
$ws.Cells.Item(8,4).Value = $null
$ws.Cells.Item(8,5).Value = $null
$ws.Cells.Item(8,6).Value = $null
$ws.Cells.Item(8,7).Value = $null
$ws.Cells.Item(8,8).Value = "Not applicable"
$ws.Cells.Item(8,9).Value = $null
$ws.Cells.Item(8,10).Value = $null
$ws.Cells.Item(8,11).Value = $null
$ws.Cells.Item(8,12).Value = $null
$ws.Cells.Item(8,13).Value = $null
$ws.Cells.Item(8,14).Value = $null
$ws.Cells.Item(8,15).Value = $null
$ws.Cells.Item(8,16).Value = $null
$ws.Cells.Item(8,17).Value = $null

# Row 9
$ws.Cells.Item(9,1).Value = 11
$ws.Cells.Item(9,2).Value = 263
$ws.Cells.Item(9,3).Value = "Chinook salmon"
$ws.Cells.Item(9,4).Value = "Spring"
$ws.Cells.Item(9,5).Value = "Juvenile"
$ws.Cells.Item(9,6).Value = "Natural"
$ws.Cells.Item(9,7).Value = "Parrot-Phelan RST"
$ws.Cells.Item(9,8).Value = "Parrott-Phelan e-test release site"
$ws.Cells.Item(9,9).Value = "n/a"
$ws.Cells.Item(9,10).Value = 249
$ws.Cells.Item(9,11).Value = 44210.5213657407
$ws.Cells.Item(9,12).Value = 6
$ws.Cells.Item(9,13).Value = "Pigment / dye"
$ws.Cells.Item(9,14).Value = "Brown"
$ws.Cells.Item(9,15).Value = "Whole body"
$ws.Cells.Item(9,16).Value = $null
$ws.Cells.Item(9,17).Value = "Yes"

# Row 10
$ws.Cells.Item(10,1).Value = 11
$ws.Cells.Item(10,2).Value = 264
$ws.Cells.Item(10,3).Value = "Chinook salmon"
$ws.Cells.Item(10,4).Value = "Spring"
$ws.Cells.Item(10,5).Value = "Juvenile"
$ws.Cells.Item(10,6).Value = "Natural"
$ws.Cells.Item(10,7).Value = "Parrot-Phelan RST"
$ws.Cells.Item(10,8).Value = "Parrott-Phelan e-test release site"
$ws.Cells.Item(10,9).Value = "n/a"
$ws.Cells.Item(10,10).Value = 109
$ws.Cells.Item(10,11).Value = 44237.5002893519
$ws.Cells.Item(10,12).Value = 7
$ws.Cells.Item(10,13).Value = "Pigment / dye"
$ws.Cells.Item(10,14).Value = "Brown"
$ws.Cells.Item(10,15).Value = "Whole body"
$ws.Cells.Item(10,16).Value = $null
$ws.Cells.Item(10,17).Value = "Yes"

# Row 11
$ws.Cells.Item(11,1).Value = 11
$ws.Cells.Item(11,2).Value = 265
$ws.Cells.Item(11,3).Value = "Chinook salmon"
$ws.Cells.Item(11,4).Value = "Spring"
$ws.Cells.Item(11,5).Value = "Juvenile"
$ws.Cells.Item(11,6).Value = "Natural"
$ws.Cells.Item(11,7).Value = "Parrot-Phelan RST"
$ws.Cells.Item(11,8).Value = "Parrott-Phelan e-test release site"
$ws.Cells.Item(11,9).Value = "n/a"
$ws.Cells.Item(11,10).Value = 349
$ws.Cells.Item(11,11).Value = 44251.5209143519
$ws.Cells.Item(11,12).Value = 7
$ws.Cells.Item(11,13).Value = "Pigment / dye"
$ws.Cells.Item(11,14).Value = "Brown"
$ws.Cells.Item(11,15).Value = "Whole body"
$ws.Cells.Item(11,16).Value = $null
$ws.Cells.Item(11,17).Value = "Yes"

# Row 12
$ws.Cells.Item(12,1).Value = 11
$ws.Cells.Item(12,2).Value = 266
$ws.Cells.Item(12,3).Value = "Chinook salmon"
$ws.Cells.Item(12,4).Value = "Spring"
$ws.Cells.Item(12,5).Value = "Juvenile"
$ws.Cells.Item(12,6).Value = "Natural"
$ws.Cells.Item(12,7).Value = "Parrot-Phelan RST"
$ws.Cells.Item(12,8).Value = "Parrott-Phelan e-test release site"
$ws.Cells.Item(12,9).Value = "n/a"
$ws.Cells.Item(12,10).Value = 197
$ws.Cells.Item(12,11).Value = 44265.5004166667
$ws.Cells.Item(12,12).Value = 7
$ws.Cells.Item(12,13).Value = "Pigment / dye"
$ws.Cells.Item(12,14).Value = "Brown"
$ws.Cells.Item(12,15).Value = "Whole body"
$ws.Cells.Item(12,16).Value = $null
$ws.Cells.Item(12,17).Value = "Yes"

# Row 13
$ws.Cells.Item(13,1).Value = 11
$ws.Cells.Item(13,2).Value = 267
$ws.Cells.Item(13,3).Value = "Chinook salmon"
$ws.Cells.Item(13,4).Value = "Spring"
$ws.Cells.Item(13,5).Value = "Juvenile"
$ws.Cells.Item(13,6).Value = "Natural"
$ws.Cells.Item(13,7).Value = "Parrot-Phelan RST"
$ws.Cells.Item(13,8).Value = "Parrott-Phelan e-test release site"
$ws.Cells.Item(13,9).Value = "n/a"
$ws.Cells.Item(13,10).Value = 160
$ws.Cells.Item(13,11).Value = 44294.5418171296
$ws.Cells.Item(13,12).Value = 7
$ws.Cells.Item(13,13).Value = "Pigment / dye"
$ws.Cells.Item(13,14).Value = "Brown"
$ws.Cells.Item(13,15).Value = "Whole body"
$ws.Cells.Item(13,16).Value = $null
$ws.Cells.Item(13,17).Value = "Yes"

# Row 14
$ws.Cells.Item(14,1).Value = 11
$ws.Cells.Item(14,2).Value = 268
$ws.Cells.Item(14,3).Value = "Chinook salmon"
$ws.Cells.Item(14,4).Value = "Spring"
$ws.Cells.Item(14,5).Value = "Juvenile"
$ws.Cells.Item(14,6).Value = "Natural"
$ws.Cells.Item(14,7).Value = "Parrot-Phelan RST"
$ws.Cells.Item(14,8).Value = "Parrott-Phelan e-test release site"
$ws.Cells.Item(14,9).Value = "n/a"
$ws.Cells.Item(14,10).Value = 500
$ws.Cells.Item(14,11).Value = 44301.5418865741
$ws.Cells.Item(14,12).Value = 7
$ws.Cells.Item(14,13).Value = "Pigment / dye"
$ws.Cells.Item(14,14).Value = "Brown"
$ws.Cells.Item(14,15).Value = "Whole body"
$ws.Cells.Item(14,16).Value = $null
$ws.Cells.Item(14,17).Value = "Yes"

# Row 15
$ws.Cells.Item(15,1).Value = 11
$ws.Cells.Item(15,2).Value = 269
$ws.Cells.Item(15,3).Value = "Chinook salmon"
$ws.Cells.Item(15,4).Value = "Spring"
$ws.Cells.Item(15,5).Value = "Juvenile"
$ws.Cells.Item(15,6).Value = "Natural"
$ws.Cells.Item(15,7).Value = "Parrot-Phelan RST"
$ws.Cells.Item(15,8).Value = "Parrott-Phelan e-test release site"
$ws.Cells.Item(15,9).Value = "n/a"
$ws.Cells.Item(15,10).Value = 250
$ws.Cells.Item(15,11).Value = 44600.5106365741
$ws.Cells.Item(15,12).Value = 7
$ws.Cells.Item(15,13).Value = "Pigment / dye"
$ws.Cells.Item(15,14).Value = "Brown"
$ws.Cells.Item(15,15).Value = "Whole body"
$ws.Cells.Item(15,16).Value = $null
$ws.Cells.Item(15,17).Value = "Yes"

# Row 16
$ws.Cells.Item(16,1).Value = 11
$ws.Cells.Item(16,2).Value = 270
$ws.Cells.Item(16,3).Value = "Chinook salmon"
$ws.Cells.Item(16,4).Value = "Spring"
$ws.Cells.Item(16,5).Value = "Juvenile"
$ws.Cells.Item(16,6).Value = "Natural"
$ws.Cells.Item(16,7).Value = "Parrot-Phelan RST"
$ws.Cells.Item(16,8).Value = "Parrott-Phelan e-test release site"
$ws.Cells.Item(16,9).Value = "n/a"
$ws.Cells.Item(16,10).Value = 249
$ws.Cells.Item(16,11).Value = 44607.5523032407
$ws.Cells.Item(16,12).Value = 7
$ws.Cells.Item(16,13).Value = "Pigment / dye"
$ws.Cells.Item(16,14).Value = "Brown"
$ws.Cells.Item(16,15).Value = "Whole body"
$ws.Cells.Item(16,16).Value = $null
$ws.Cells.Item(16,17).Value = "Yes"

# Row 17
$ws.Cells.Item(17,1).Value = 11
$ws.Cells.Item(17,2).Value = 271
$ws.Cells.Item(17,3).Value = "Chinook salmon"
$ws.Cells.Item(17,4).Value = "Spring"
$ws.Cells.Item(17,5).Value = "Juvenile"
$ws.Cells.Item(17,6).Value = "Natural"
$ws.Cells.Item(17,7).Value = $null
$ws.Cells.Item(17,8).Value = "Parrott-Phelan e-test release site"
$ws.Cells.Item(17,9).Value = "n/a"
$ws.Cells.Item(17,10).Value = 185
$ws.Cells.Item(17,11).Value = 44628.555775463
$ws.Cells.Item(17,12).Value = 7
$ws.Cells.Item(17,13).Value = "Pigment / dye"
$ws.Cells.Item(17,14).Value = "Brown"
$ws.Cells.Item(17,15).Value = "Whole body"
$ws.Cells.Item(17,16).Value = $null
$ws.Cells.Item(17,17).Value = "Yes"

# Row 18
$ws.Cells.Item(18,1).Value = $null
$ws.Cells.Item(18,2).Value = $null
$ws.Cells.Item(18,3).Value = $null
$ws.Cells.Item(18,4).Value = $null
$ws.Cells.Item(18,5).Value = $null
$ws.Cells.Item(18,6).Value = $null
$ws.Cells.Item(18,7).Value = $null
$ws.Cells.Item(18,8).Value = $null
$ws.Cells.Item(18,9).Value = $null
$ws.Cells.Item(18,10).Value = $null
$ws.Cells.Item(18,11).Value = $null
$ws.Cells.Item(18,12).Value = $null
$ws.Cells.Item(18,13).Value = $null
$ws.Cells.Item(18,14).Value = $null
$ws.Cells.Item(18,15).Value = $null
$ws.Cells.Item(18,16).Value = $null
$ws.Cells.Item(18,17).Value = $null

# Row 19
$ws.Cells.Item(19,1).Value = $null
$ws.Cells.Item(19,2).Value = $null
$ws.Cells.Item(19,3).Value = $null
$ws.Cells.Item(19,4).Value = $null
$ws.Cells.Item(19,5).Value = $null
$ws.Cells.Item(19,6).Value = $null
$ws.Cells.Item(19,7).Value = $null
$ws.Cells.Item(19,8).Value = $null
$ws.Cells.Item(19,9).Value = $null
$ws.Cells.Item(19,10).Value = $null
$ws.Cells.Item(19,11).Value = $null
$ws.Cells.Item(19,12).Value = $null
$ws.Cells.Item(19,13).Value = $null
$ws.Cells.Item(19,14).Value = $null
$ws.Cells.Item(19,15).Value = $null
$ws.Cells.Item(19,16).Value = $null
$ws.Cells.Item(19,17).Value = $null

# Row 20
$ws.Cells.Item(20,1).Value = $null
$ws.Cells.Item(20,2).Value = $null
$ws.Cells.Item(20,3).Value = $null
$ws.Cells.Item(20,4).Value = $null
$ws.Cells.Item(20,5).Value = $null
$ws.Cells.Item(20,6).Value = $null
$ws.Cells.Item(20,7).Value = $null
$ws.Cells.Item(20,8).Value = $null
$ws.Cells.Item(20,9).Value = $null
$ws.Cells.Item(20,10).Value = $null
$ws.Cells.Item(20,11).Value = $null
$ws.Cells.Item(20,12).Value = $null
$ws.Cells.Item(20,13).Value = $null
$ws.Cells.Item(20,14).Value = $null
$ws.Cells.Item(20,15).Value = $null
$ws.Cells.Item(20,16).Value = $null
$ws.Cells.Item(20,17).Value = $null

# Row 21
$ws.Cells.Item(21,1).Value = $null
$ws.Cells.Item(21,2).Value = $null
$ws.Cells.Item(21,3).Value = $null
$ws.Cells.Item(21,4).Value = $null
$ws.Cells.Item(21,5).Value = $null
$ws.Cells.Item(21,6).Value = $null
$ws.Cells.Item(21,7).Value = $null
$ws.Cells.Item(21,8).Value = $null
$ws.Cells.Item(21,9).Value = $null
$ws.Cells.Item(21,10).Value = $null
$ws.Cells.Item(21,11).Value = $null
$ws.Cells.Item(21,12).Value = $null
$ws.Cells.Item(21,13).Value = $null
$ws.Cells.Item(21,14).Value = $null
$ws.Cells.Item(21,15).Value = $null
$ws.Cells.Item(21,16).Value = $null
$ws.Cells.Item(21,17).Value = $null

# Update defined name range
$n = $wb.Names.Item("Release_EDI")
$n.RefersTo = "=Release_EDI!`$A`$1:`$Q`$17"
